# Slide 10, "Content Placeholder 2" shape: the last bullet paragraph currently
# reads across three runs:
#   "Does not modify existing STAMP (which is for DM) procedure as different UDP destination "
#   "Port2"   (bold, blue)
#   " is used for LM"
# It should become a single run (matching the first run's formatting):
#   "Does not modify existing STAMP (which is for DM) procedure as different destination UDP is used for LM"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(3)               # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(7, 1)           # the target bullet paragraph

# Replace the whole paragraph's characters (not just the .Text property) so
# the three existing runs collapse into a single run carrying the first
# run's formatting (non-bold, tx2 scheme color).
$full = $para.Characters(1, $para.Length)
$full.Text = "Does not modify existing STAMP (which is for DM) procedure as different destination UDP is used for LM"
